$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 532 (shifts existing rows 532:636 down to 533:637)
$ws.Rows.Item(532).Insert()

# Populate the newly inserted row with the new observation
$ws.Range("A532").Value = 8
$ws.Range("B532").Value = "Terminal La Palmera de La Serena"
$ws.Range("C532").Value = "Coquimbo"
$ws.Range("D532").Value = 45275
$ws.Range("E532").Value = 4
$ws.Range("F532").Value = 100112032
$ws.Range("G532").Value = "Zapallo italiano"
$ws.Range("H532").Value = "Sin especificar"
$ws.Range("I532").Value = "Primera"
$ws.Range("J532").Value = 400
$ws.Range("K532").Value = 8000
$ws.Range("L532").Value = 9000
$ws.Range("M532").Value = 8500
$ws.Range("N532").Value = "`$/caja 60 unidades"
$ws.Range("O532").Value = "Provincia del Elquí"
$ws.Range("P532").Value = 142
$ws.Range("Q532").Value = 60
$ws.Range("R532").Value = "Hortaliza"
